$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(", and report to jail on September 23, 2022, at 7:00 p.m.", $false, $false, `
              $false, $false, $false, $true, 1, $false, ".", 2)
